$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "dhcp"
# ---------------------------------------------------------------------------
$dhcp = $wb.Worksheets.Item("dhcp")

# Rename existing "Printer" host row to "FinancePrinter3"
$dhcp.Range("A4").Value = "FinancePrinter3"

# Add the new "AdminProjector1" host row
$dhcp.Range("A5").Value = "AdminProjector1"
$dhcp.Range("B5").Value = "172.24.161.7"
$dhcp.Range("C5").Value = "None"
$dhcp.Range("D5").Value = "None"
$dhcp.Range("E5").Value = "255.255.255.0"
$dhcp.Range("F5").Value = "172.24.161.1"
$dhcp.Range("G5").Value = "172.24.161.7"
$dhcp.Range("H5").Value = "0144.5566.DDEE.FF"
$dhcp.Range("I5").Value = "static"

# Widen column A to fit the new, longer host names
$dhcp.Columns.Item(1).ColumnWidth = 16.5

# ---------------------------------------------------------------------------
# Sheet 2: "vlan"
# ---------------------------------------------------------------------------
$vlan = $wb.Worksheets.Item("vlan")

# New "description" column
$vlan.Range("F1").Value = "description"
$vlan.Range("F2").Value = "Connection_to_Core_1"
$vlan.Range("F3").Value = "Connection_to_Core_2"
$vlan.Range("F4").Value = "Connection_to_Data"
$vlan.Range("F5").Value = "Connection_to_Voice"
$vlan.Range("F6").Value = "None"
$vlan.Range("F7").Value = "None"
$vlan.Range("F8").Value = "None"
$vlan.Range("F9").Value = "None"

# New "Management" vlan row
$vlan.Range("A10").Value = "Management"
$vlan.Range("B10").Value = 900
$vlan.Range("C10").Value = "UPLINK"
$vlan.Range("D10").Value = "172.20.160.1"
$vlan.Range("E10").Value = "255.255.255.0"
$vlan.Range("F10").Value = "Connection_to_Management"

# Column width adjustments
$vlan.Columns.Item(2).ColumnWidth = 10.166666666666666
$vlan.Columns.Item(6).ColumnWidth = 31.166666666666668

# ---------------------------------------------------------------------------
# Sheet 3: "intf_access" (new sheet, appended after "vlan")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$intf = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$intf.Name = "intf_access"

$intf.Range("A1").Value = "description"
$intf.Range("B1").Value = "access_vlan"
$intf.Range("C1").Value = "port_start"
$intf.Range("D1").Value = "port_end"

$intf.Range("A2").Value = "Connection_to_Users_&_Phones"
$intf.Range("B2").Value = 20
$intf.Range("C2").Value = "FastEthernet0/1"
$intf.Range("D2").Value = "FastEthernet0/2"

$intf.Range("A3").Value = "Connection_to_CCTV"
$intf.Range("B3").Value = 200
$intf.Range("C3").Value = "FastEthernet0/3"
$intf.Range("D3").Value = "FastEthernet0/5"

$intf.Range("A4").Value = "Connection_to_AP_MGT"
$intf.Range("B4").Value = 600
$intf.Range("C4").Value = "FastEthernet0/7"
$intf.Range("D4").Value = "None"

$intf.Columns.Item(1).ColumnWidth = 29.666666666666668
$intf.Columns.Item(2).ColumnWidth = 13.666666666666666
$intf.Columns.Item(3).ColumnWidth = 21.166666666666668
$intf.Columns.Item(4).ColumnWidth = 19.666666666666668

# ---------------------------------------------------------------------------
# Final selections per sheet (also drives which sheet/tab ends up active)
# ---------------------------------------------------------------------------
$dhcp.Activate()
$dhcp.Range("H8").Select() | Out-Null

$vlan.Activate()
$vlan.Range("C14").Select() | Out-Null

$intf.Activate()
$intf.Range("C13").Select() | Out-Null
